$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = -8
    4  = -5
    5  = 2
    6  = -5
    9  = -4
    11 = -6
    12 = -1
    13 = -3
    16 = 8
    19 = 4
    24 = 0
    25 = -3
    27 = 2
    30 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
